# "exchanged term from genbank 'mapping'"
#
# The GENBANK_GENOME reference sheet used a leftover NFDI4PSO term
# ("Parameter [Processed data file name]" / NFDI4PSO:0000028) to describe
# the filename of the processed/derived data. That term is retired in favor
# of the DataPLANT "Derived Data File" term, which now carries the same
# filename guidance. The row that used to hold
# "Parameter [Processed data file format]" (DPBO:0000027) moves up one row
# to take the old row's place, and a new "Derived Data File" row is appended
# at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GENBANK_GENOME")

# Row 15: was "Parameter [Processed data file name]" (NFDI4PSO) -> becomes
# "Parameter [Processed data file format]" (DPBO:0000027), i.e. what used
# to live on row 16.
$ws.Range("A15").Value = "Parameter [Processed data file format]"
$ws.Range("B15").Value = "DPBO:0000027"
$ws.Range("C15").Value = "DPBO"
$ws.Range("D15").Value = "http://purl.obolibrary.org/obo/DPBO_0000027"
$ws.Range("G15").Value = ""
$ws.Range("H15").Value = ""
$ws.Range("I15").Value = "n"

# Row 16: becomes the new "Derived Data File" row, carrying the filename
# notes/instruction/requirement that used to live on the removed row 15.
$ws.Range("A16").Value = "Derived Data File"
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = ""
$ws.Range("G16").Value = "filename"
$ws.Range("H16").Value = "Each genome must be in 1 or 2 files. You can concatenate multiple fasta files into a single file that can be submitted or used as input for tbl2asn.`n"
$ws.Range("I16").Value = "m"

# Rebuild the term-accession hyperlinks in column D to match the new
# row layout (the NFDI4PSO_0000028 link on the old row 15 is gone; the
# DPBO_0000027 link now sits on row 15 instead of row 16).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D4"), "http://purl.obolibrary.org/obo/NCIT_C175889") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "http://purl.obolibrary.org/obo/DPBO_0000023") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "http://purl.obolibrary.org/obo/DPBO_0000024") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "http://purl.obolibrary.org/obo/DPBO_0000025") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), "http://purl.obolibrary.org/obo/DPBO_0000040") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D10"), "http://purl.obolibrary.org/obo/DPBO_0000060") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D13"), "http://purl.obolibrary.org/obo/DPBO_0000061") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D14"), "http://purl.obolibrary.org/obo/DPBO_0000026") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D15"), "http://purl.obolibrary.org/obo/DPBO_0000027") | Out-Null
